# Updated symbol list on Wed Dec 14 05:12:08 UTC 2022 with GitHub Actions
#
# Refreshes the cryptos.xlsx price table: for each coin row (2-51) the
# "Hora" column (G) advances from 4 to 5, most rows get a refreshed
# "Price" (D), and two rows have their "Volume(1h)" label (E) gain/lose
# the "Worstin24h" suffix as the ranking shuffles.
#
# Price/Hora values are numeric-looking text (e.g. "0.09300", "0.001090")
# where trailing zeros are significant, so they are written with a
# leading apostrophe to force Excel to keep them as text instead of
# normalizing them into floating point numbers. The cell style is then
# reset to "Normal" so no stray quote-prefix formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

$rowUpdates = @(
    @{ Row = 2; D = "275.53"; E = $null }
    @{ Row = 3; D = "23.12"; E = $null }
    @{ Row = 4; D = "6.451"; E = $null }
    @{ Row = 5; D = "0.06273"; E = $null }
    @{ Row = 6; D = "3.654"; E = $null }
    @{ Row = 7; D = "6.668"; E = $null }
    @{ Row = 8; D = "1.403"; E = $null }
    @{ Row = 9; D = "0.8309"; E = $null }
    @{ Row = 10; D = $null; E = $null }
    @{ Row = 11; D = "0.1626"; E = $null }
    @{ Row = 12; D = "0.08306"; E = $null }
    @{ Row = 13; D = $null; E = $null }
    @{ Row = 14; D = "0.03123"; E = $null }
    @{ Row = 15; D = "0.09300"; E = $null }
    @{ Row = 16; D = "3.844"; E = $null }
    @{ Row = 17; D = "0.001647"; E = $null }
    @{ Row = 18; D = "0.04783"; E = $null }
    @{ Row = 19; D = "0.006293"; E = $null }
    @{ Row = 20; D = "0.005681"; E = "19HotbitTokenHTBWorstin24h" }
    @{ Row = 21; D = "0.001090"; E = $null }
    @{ Row = 22; D = $null; E = $null }
    @{ Row = 23; D = "3.715"; E = $null }
    @{ Row = 24; D = "2.323"; E = $null }
    @{ Row = 25; D = "0.3346"; E = $null }
    @{ Row = 26; D = "0.1240"; E = $null }
    @{ Row = 27; D = "0.0002680"; E = $null }
    @{ Row = 28; D = $null; E = $null }
    @{ Row = 29; D = $null; E = $null }
    @{ Row = 30; D = $null; E = $null }
    @{ Row = 31; D = $null; E = $null }
    @{ Row = 32; D = $null; E = $null }
    @{ Row = 33; D = $null; E = $null }
    @{ Row = 34; D = $null; E = $null }
    @{ Row = 35; D = $null; E = $null }
    @{ Row = 36; D = $null; E = $null }
    @{ Row = 37; D = $null; E = $null }
    @{ Row = 38; D = $null; E = $null }
    @{ Row = 39; D = $null; E = $null }
    @{ Row = 40; D = "0.04718"; E = $null }
    @{ Row = 41; D = "0.007048"; E = $null }
    @{ Row = 42; D = "0.1162"; E = $null }
    @{ Row = 43; D = $null; E = "42CEJICEJI" }
    @{ Row = 44; D = "0.01219"; E = $null }
    @{ Row = 45; D = "0.00006253"; E = $null }
    @{ Row = 46; D = $null; E = $null }
    @{ Row = 47; D = $null; E = $null }
    @{ Row = 48; D = "0.7965"; E = $null }
    @{ Row = 49; D = "0.03444"; E = $null }
    @{ Row = 50; D = "0.00002300"; E = $null }
    @{ Row = 51; D = $null; E = $null }
)

foreach ($update in $rowUpdates) {
    $row = $update.Row

    if ($null -ne $update.D) {
        Set-TextValue $ws.Cells.Item($row, 4) $update.D
    }

    if ($null -ne $update.E) {
        $ws.Cells.Item($row, 5).Value = $update.E
    }

    Set-TextValue $ws.Cells.Item($row, 7) "5"
}
